$d = $word.ActiveDocument

# --- Move the "_GoBack" bookmark ---
# It currently sits between "RA DESEMPEÑAR EL PUEST" and "O DE " (inside the
# "...SUFICIENTES PARA DESEMPEÑAR EL PUESTO DE ANALISTA." sentence). The edit
# removes it from there (that text becomes one contiguous, unbroken run) and
# re-creates it further up the document, right after the new
# "...CELEBRA EL DIA 8 DE ENERO DEL 2019," text (before " POR UNA ").
$hadGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hadGoBack) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 1: opening recital sentence ---
$rng1 = $d.Content
$rng1.Find.Execute(
    "CONTRATO INDIVIDUAL DE TRABAJO QUE CELEBRARAN, POR UNA ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CONTRATO INDIVIDUAL DE TRABAJO QUE CELEBRA EL DIA 8 DE ENERO DEL 2019, POR UNA ",
    2
) | Out-Null

# Re-create _GoBack right before " POR UNA " (i.e. right after "2019,").
$rngBm = $d.Content
$rngBm.Find.Execute("2019, POR UNA ") | Out-Null
$bmPoint = $rngBm.Start + 5   # length of "2019," = 5 chars
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Change 3 & 4: employment dates sentence ---
$rng3 = $d.Content
$rng3.Find.Execute(
    "ENERO HASTA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ENERO DEL 2019 HASTA",
    2
) | Out-Null

$rng4 = $d.Content
$rng4.Find.Execute(
    "EL 21 DE AGOSTO.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "EL 21 DE AGOSTO DEL 2019.",
    2
) | Out-Null

Write-Host "Done."
